$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for column F
$ws.Range("F1").Value = "IsActive"

# Add boolean values for rows 2-4
$ws.Range("F2").Value = $true
$ws.Range("F3").Value = $false
$ws.Range("F4").Value = $true

# Apply border style to F2:F4 (matching existing E column formatting look)
$ws.Range("F1:F4").Borders.LineStyle = 1

# Update selection to I5
$ws.Range("I5").Select()
